# Generate Report for Handoff
# The dfc7aa25-1809-4c6f-a4f4-dfce41147197.md file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with new
# handoff timestamps and a stale-handback warning recorded in the
# per-language sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is the dfc7aa25... file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-02 06:56:29"

# ---- zh-cn sheet: row 3 is the dfc7aa25... file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-09-02 06:56:24"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/040529efba701779b5c2a6de3228b932212551b2/e2e/dfc7aa25-1809-4c6f-a4f4-dfce41147197.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7041629eacfde58df0aa6b5e4bf8797c6fd258a/e2e/dfc7aa25-1809-4c6f-a4f4-dfce41147197.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.2

# ---- de-de sheet: row 3 is the dfc7aa25... file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-09-02 06:56:29"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/040529efba701779b5c2a6de3228b932212551b2/e2e/dfc7aa25-1809-4c6f-a4f4-dfce41147197.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7041629eacfde58df0aa6b5e4bf8797c6fd258a/e2e/dfc7aa25-1809-4c6f-a4f4-dfce41147197.md."
$dede.Columns.Item(16).ColumnWidth = 39.2
